$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.6498683146729446
$ws.Range("D2").Value = 0.5225057998191049

$ws.Range("C3").Value = -0.3925152442352722
$ws.Range("D3").Value = 0.6984546676101222

$ws.Range("C4").Value = 0.1505171351633422
$ws.Range("D4").Value = 0.8817281538429809

$ws.Range("C5").Value = 2.678961337786918
$ws.Range("D5").Value = 0.01371002736325333

$ws.Range("C6").Value = -1.071356376449245
$ws.Range("D6").Value = 0.29562219903482

$ws.Range("C7").Value = -0.5505560020292669
$ws.Range("D7").Value = 0.5874855913056884

$ws.Range("C8").Value = 1.455511616841284
$ws.Range("D8").Value = 0.1596500353693517

$ws.Range("C9").Value = 0.6233217515520759
$ws.Range("D9").Value = 0.5394821841408546

$ws.Range("C10").Value = 2.10804248464454
$ws.Range("D10").Value = 0.04664727541731994

$ws.Range("C11").Value = 1.926996697664108
$ws.Range("D11").Value = 0.06699313630424397
$ws.Range("G11").Value = "No"
